$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading-percent results for rows 2-25 (index A = 0..23),
# columns B..N (E, M, O remain 0; column A unchanged) - "case with 380 kV done"
$data = @(
@(15.45322887519812,5.454524623665281,7.630446730495924,0,39.93366610473762,46.51529974670029,18.79429747892872,29.44411109670961,10.34789756718635,12.68037904236986,11.92238557940252,0,21.48115504429802),
@(15.27558165316948,5.337569766701939,7.618021024537695,0,39.99408338023418,46.57514058803697,18.84071145555603,29.52384162097812,10.3671935495811,12.55726033829397,11.92340868948161,0,21.54395961560254),
@(15.1686462358763,5.263456402650477,7.611346060578951,0,40.03903686003868,46.6228518969219,18.87196216312665,29.57757617779722,10.37984214586248,12.48356421094219,11.92568912131474,0,21.58433364710512),
@(15.12565646217693,5.232696202505588,7.608868052404941,0,40.05932945374339,46.64504784276862,18.88538904482472,29.6006744829685,10.38519836486402,12.4540409380973,11.9270350028209,0,21.60124326701005),
@(15.11855492348335,5.227555394694581,7.608471261167376,0,40.06281816055807,46.64889958333892,18.88766035696589,29.60458245276768,10.38609996321341,12.44917021791684,11.92728367646396,0,21.60407872855282),
@(15.16806401897389,5.263043792655103,7.611311658403443,0,40.03930254473852,46.62314009900214,18.87214044121608,29.57788282784834,10.37991356397184,12.4831639507393,11.92570558421178,0,21.58455984458455),
@(15.39156150043876,5.414688473041633,7.625965535727442,0,39.95286617103161,46.53365358992788,18.8097296969162,29.47060965587741,10.35438488919917,12.63755234498302,11.92239613127449,0,21.50243480752829),
@(15.84463066316457,5.692952749089073,7.662178585187165,0,39.84578098516543,46.44538620377131,18.70918665017162,29.29820928228941,10.3106569463916,12.95397891837028,11.92896303287568,0,21.35570350203348),
@(16.18354486518517,5.884709795468531,7.693206086737211,0,39.80523310126929,46.43389509631222,18.64864174477625,29.19473880186336,10.28236350483758,13.19287651846865,11.94167048058099,0,21.25654175219159),
@(16.33840609202737,5.968997466650637,7.70824973498604,0,39.79507161139819,46.44027131956187,18.62399264218469,29.15271303126837,10.27031863892645,13.30253689638535,11.94914403328444,0,21.21328874620404),
@(16.39709460257484,6.000477131912215,7.714077062135967,0,39.79241434146133,46.44435368304551,18.61507475667987,29.13752494445847,10.26587588570863,13.34416924162748,11.95221563454794,0,21.19717553459407),
@(16.38445377219029,5.993717139347803,7.712816282139532,0,39.79293369212079,46.4434003213163,18.6169768700953,29.1407636606568,10.26682745357043,13.33519880020399,11.95154339998347,0,21.20063400135109),
@(16.34323381063578,5.971596182198365,7.708726552285309,0,39.7948291387964,46.44057376110654,18.62325062033818,29.1514489398845,10.26995076076943,13.30596008428616,11.94939190656577,0,21.21195778285782),
@(16.31798981926614,5.957988961066319,7.706238395315604,0,39.79614518774422,46.43905956784455,18.62714767894028,29.15808857549209,10.27187928097312,13.28806333312342,11.94810545006063,0,21.21892850219851),
@(16.17343360881123,5.879140858866057,7.692241392822963,0,39.80606383882993,46.43371185279963,18.6503108533056,29.19758685210487,10.28316725167073,13.18572675018944,11.94121597167651,0,21.25940569000744),
@(16.08489120592432,5.83000514963253,7.683890567031512,0,39.81427012107248,46.43340301686601,18.66526183101089,29.22311019161506,10.29030332178484,13.1231736207081,11.93742177451854,0,21.28471172016972),
@(16.03403178599934,5.801467342959341,7.679175028708054,0,39.81976995245689,46.43431769061416,18.6741335397581,29.23826520597898,10.29448556754701,13.08728949776262,11.93539891646973,0,21.29944185911165),
@(16.01682493655358,5.791757962243516,7.677593568642357,0,39.8217660385921,46.43481499107385,18.6771841102078,29.24347792858097,10.29591497278884,13.07515710692328,11.93474145527336,0,21.30445928496137),
@(16.09431006207367,5.835264432463488,7.68477047789369,0,39.81331584457772,46.43332284467662,18.66364209032068,29.22034405568774,10.28953562923128,13.12982295456298,11.9378091827284,0,21.28199976938245),
@(16.3553402834349,5.978105650697095,7.709924284315233,0,39.79424009322793,46.44135873824028,18.62139657118706,29.14829069961022,10.26903016018259,13.31454559754766,11.95001731316843,0,21.20862451012655),
@(16.52618245717862,6.068900109374594,7.727123661267282,0,39.78871241604373,46.45633079185602,18.59621264069293,29.10543266151522,10.25631845586131,13.43587622071336,11.95940286714791,0,21.16221808072698),
@(16.43499570065174,6.020680321045988,7.717875513207886,0,39.79102802180155,46.44745110501117,18.6094317541788,29.12791919874471,10.26303994450285,13.37107611714083,11.95426555969458,0,21.18684475117875),
@(16.09005165364524,5.8328876098588,7.684372403776293,0,39.8137448372558,46.43335568822148,18.66437351475521,29.22159312654177,10.28988245498334,13.12681654349428,11.93763354170068,0,21.28322527728028),
@(15.72078624980368,5.619832400617722,7.651595081163241,0,39.8680581832582,46.45990372055964,18.734046782142,29.34077874288607,10.32181132531801,12.86710961940289,11.9257961068565,0,21.39387449700152)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowValues[$j]
    }
}
